# Stundenliste-MichaelSchneider.xlsx update
# - add 3 new logged entries (rows 40-42) with date / hours / running total / activity
# - move the view's scroll position / active selection down to the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
# Row 40: 06.09.2021, 1h, "internes meeting"
$ws.Range("A40").Value = 44445
$ws.Range("B40").Value = 1
$ws.Range("D40").Value = "internes meeting"

# Row 41: 07.09.2021, 5h, "Dokumentation, tests"
$ws.Range("A41").Value = 44446
$ws.Range("B41").Value = 5
$ws.Range("D41").Value = "Dokumentation, tests"

# Row 42: 09.09.2021, 2h, "Klassendiagramm gecheckt und updated"
$ws.Range("A42").Value = 44448
$ws.Range("B42").Value = 2
$ws.Range("D42").Value = "Klassendiagramm gecheckt und updated"

# Running-total formula (same pattern as the existing C11:C39 shared formula)
$ws.Range("C40:C42").Formula = "=C39+B40"

# Match the style used by the other "Taetigkeit" cells in that column (fontId 3 /
# cellXf 9) by copying the formatting from the cell directly above (D39) instead
# of assigning a brand new style object.
$ws.Range("D39").Copy() | Out-Null
$ws.Range("D40:D42").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- View / selection ------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1
$ws.Range("D43").Select() | Out-Null
